$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3594.1904
$ws.Range("J17").Value = 3553.805
$ws.Range("L17").Value = 10661.415
$ws.Range("N17").Value = -10997.415

$ws.Range("H32").Value = 2522
$ws.Range("J32").Value = 2963.1667
$ws.Range("L32").Value = 2963.1667
$ws.Range("N32").Value = -3615.1667

$ws.Range("H43").Value = 19444.111
$ws.Range("J43").Value = 19856.715
$ws.Range("L43").Value = 19856.715
$ws.Range("N43").Value = -19994.715

$ws.Range("H113").Value = 2850.7334
$ws.Range("I113").Value = 2971.5
$ws.Range("J113").Value = 2770.2222
$ws.Range("K113").Value = 2971.5
$ws.Range("L113").Value = 2770.2222
$ws.Range("M113").Value = 282.5
$ws.Range("N113").Value = -9278.2222

$ws.Range("H116").Value = 3333.111
$ws.Range("J116").Value = 3499.5
$ws.Range("L116").Value = 3499.5
$ws.Range("N116").Value = -10383.5

$ws.Range("H132").Value = 2145.762
$ws.Range("I132").Value = 1848
$ws.Range("J132").Value = 3932.3333
$ws.Range("K132").Value = 5544
$ws.Range("L132").Value = 11796.9999
$ws.Range("M132").Value = -3014
$ws.Range("N132").Value = -16856.9999

$ws.Range("H137").Value = 29260.467
$ws.Range("I137").Value = 59754.785
$ws.Range("J137").Value = 2577.9375
$ws.Range("K137").Value = 179264.355
$ws.Range("L137").Value = 7733.8125
$ws.Range("M137").Value = -176714.355
$ws.Range("N137").Value = -12833.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20748716
$ws.Range("I32").Value = 19930064
$ws.Range("K32").Value = 19930064
$ws.Range("M32").Value = -19929777

$ws.Range("H45").Value = 4120.375
$ws.Range("I45").Value = 3780.4285
$ws.Range("K45").Value = 3780.4285
$ws.Range("M45").Value = -3403.4285

$ws.Range("H61").Value = 3626.5386
$ws.Range("I61").Value = 3397.84
$ws.Range("K61").Value = 3397.84
$ws.Range("M61").Value = -3185.84

$ws.Range("H63").Value = 6412.5
$ws.Range("J63").Value = 10325
$ws.Range("L63").Value = 10325
$ws.Range("N63").Value = -11697

$ws.Range("H66").Value = 6412.5
$ws.Range("J66").Value = 10325
$ws.Range("L66").Value = 51625
$ws.Range("N66").Value = -58489

$ws.Range("H132").Value = 259290.08
$ws.Range("I132").Value = 306008.03
$ws.Range("J132").Value = 2341.3333
$ws.Range("K132").Value = 918024.0900000001
$ws.Range("L132").Value = 7023.999899999999
$ws.Range("M132").Value = -915494.0900000001
$ws.Range("N132").Value = -12083.9999

$ws.Range("H136").Value = 3626.5386
$ws.Range("I136").Value = 3397.84
$ws.Range("K136").Value = 10193.52
$ws.Range("M136").Value = -7643.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3425.4443
$ws.Range("I86").Value = 2726.1667
$ws.Range("K86").Value = 2726.1667
$ws.Range("M86").Value = -1603.1667

$ws.Range("H89").Value = 3425.4443
$ws.Range("I89").Value = 2726.1667
$ws.Range("K89").Value = 13630.8335
$ws.Range("M89").Value = -8014.833500000001

$ws.Range("H94").Value = 583.7143
$ws.Range("I94").Value = 477.88
$ws.Range("J94").Value = 1465.6666
$ws.Range("K94").Value = 477.88
$ws.Range("L94").Value = 1465.6666
$ws.Range("M94").Value = -26.88
$ws.Range("N94").Value = -2367.6666

$ws.Range("H99").Value = 2497.647
$ws.Range("J99").Value = 3747.25
$ws.Range("L99").Value = 3747.25
$ws.Range("N99").Value = -6743.25

$ws.Range("H105").Value = 4089.9
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 4000
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2832.86
$ws.Range("I31").Value = 1267.48
$ws.Range("J31").Value = 4398.24
$ws.Range("K31").Value = 1267.48
$ws.Range("L31").Value = 4398.24
$ws.Range("M31").Value = -972.48
$ws.Range("N31").Value = -4988.24

$ws.Range("H34").Value = 2832.86
$ws.Range("I34").Value = 1267.48
$ws.Range("J34").Value = 4398.24
$ws.Range("K34").Value = 1267.48
$ws.Range("L34").Value = 4398.24
$ws.Range("M34").Value = -1065.48
$ws.Range("N34").Value = -4802.24

$ws.Range("H58").Value = 3341.3096
$ws.Range("I58").Value = 2745.2173
$ws.Range("J58").Value = 4062.8948
$ws.Range("K58").Value = 2745.2173
$ws.Range("L58").Value = 4062.8948
$ws.Range("M58").Value = -2542.2173
$ws.Range("N58").Value = -4468.8948

$ws.Range("H99").Value = 2542.611
$ws.Range("I99").Value = 2419.7273
$ws.Range("K99").Value = 2419.7273
$ws.Range("M99").Value = -921.7273

$ws.Range("H120").Value = 9666
$ws.Range("J120").Value = 9666
$ws.Range("L120").Value = 9666
$ws.Range("N120").Value = -16924

$ws.Range("H126").Value = 2542.611
$ws.Range("I126").Value = 2419.7273
$ws.Range("K126").Value = 7259.1819
$ws.Range("M126").Value = -4789.1819

$ws.Range("H132").Value = 3268.4126
$ws.Range("I132").Value = 3029.2126
$ws.Range("J132").Value = 3971.0625
$ws.Range("K132").Value = 9087.6378
$ws.Range("L132").Value = 11913.1875
$ws.Range("M132").Value = -6557.6378
$ws.Range("N132").Value = -16973.1875

$ws.Range("H136").Value = 3341.3096
$ws.Range("I136").Value = 2745.2173
$ws.Range("J136").Value = 4062.8948
$ws.Range("K136").Value = 8235.651899999999
$ws.Range("L136").Value = 12188.6844
$ws.Range("M136").Value = -5685.651899999999
$ws.Range("N136").Value = -17288.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1351.1034
$ws.Range("J107").Value = 1400.8096
$ws.Range("L107").Value = 4202.4288
$ws.Range("N107").Value = -8042.4288

$ws.Range("H113").Value = 2356.2
$ws.Range("J113").Value = 3429.6667
$ws.Range("L113").Value = 10289.0001
$ws.Range("N113").Value = -14629.0001

$ws.Range("H131").Value = 1581.8596
$ws.Range("I131").Value = 943.0769
$ws.Range("J131").Value = 1770.591
$ws.Range("K131").Value = 2829.2307
$ws.Range("L131").Value = 5311.772999999999
$ws.Range("M131").Value = 2210.7693
$ws.Range("N131").Value = -15391.773

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 38680.184
$ws.Range("I113").Value = 5743.857
$ws.Range("J113").Value = 54050.465
$ws.Range("K113").Value = 5743.857
$ws.Range("L113").Value = 54050.465
$ws.Range("M113").Value = -3573.857
$ws.Range("N113").Value = -58390.465

$ws.Range("H126").Value = 2640.4443
$ws.Range("I126").Value = 2561.0833
$ws.Range("K126").Value = 7683.249899999999
$ws.Range("M126").Value = -5213.249899999999

$ws.Range("H132").Value = 2723.3555
$ws.Range("I132").Value = 2446.4827
$ws.Range("J132").Value = 3225.1875
$ws.Range("K132").Value = 7339.4481
$ws.Range("L132").Value = 9675.5625
$ws.Range("M132").Value = -4809.4481
$ws.Range("N132").Value = -14735.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2131.45
$ws.Range("I22").Value = 1892.7273
$ws.Range("J22").Value = 2423.2222
$ws.Range("K22").Value = 1892.7273
$ws.Range("L22").Value = 2423.2222
$ws.Range("M22").Value = -1597.7273
$ws.Range("N22").Value = -3013.2222

$ws.Range("H27").Value = 2131.45
$ws.Range("I27").Value = 1892.7273
$ws.Range("J27").Value = 2423.2222
$ws.Range("K27").Value = 1892.7273
$ws.Range("L27").Value = 2423.2222
$ws.Range("M27").Value = -1785.7273
$ws.Range("N27").Value = -2637.2222

$ws.Range("H46").Value = 6521.3794
$ws.Range("J46").Value = 7885
$ws.Range("L46").Value = 7885
$ws.Range("N46").Value = -8261

$ws.Range("H55").Value = 1299.625
$ws.Range("I55").Value = 1233.4445
$ws.Range("J55").Value = 1498.1666
$ws.Range("K55").Value = 1233.4445
$ws.Range("L55").Value = 1498.1666
$ws.Range("M55").Value = -1060.4445
$ws.Range("N55").Value = -1844.1666

$ws.Range("H136").Value = 6080.1904
$ws.Range("I136").Value = 6340.4
$ws.Range("J136").Value = 5429.6665
$ws.Range("K136").Value = 19021.2
$ws.Range("L136").Value = 16288.9995
$ws.Range("M136").Value = -16471.2
$ws.Range("N136").Value = -21388.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1670369.5
$ws.Range("I4").Value = 1254165.8
$ws.Range("K4").Value = 1254165.8
$ws.Range("M4").Value = -1254052.8
